$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C3DC queries were updated to join on the renamed id columns
# (study_id / participant_id) instead of the old generic "id" columns.
$pairs = @(
  @('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"'),
  @('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'),
  @('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'),
  @('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'),
  @('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"'),
  @('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
)

foreach ($pair in $pairs) {
  $ws.Cells.Replace($pair[0], $pair[1]) | Out-Null
}

# Widen column C (StatQuery) now that the queries are longer, and let Excel
# recompute the width naturally instead of keeping the old "best fit" flag.
$ws.Columns.Item(3).ColumnWidth = 68.5

# Scroll the view down one row (topLeftCell A6 -> A7) while keeping the
# current selection on C7.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select() | Out-Null
